# 自动更新Excel文件 - 2025-11-11 23:20:55
# This script advances the "countdown" tracker sheet by one day:
#   - Column E ("剩余" / days remaining) is decremented by 1 for every
#     data row (rows 2-99), EXCEPT row 36 whose start-date cell already
#     contains a malformed 9-digit value and whose remaining count equals
#     its total, so it is left untouched (matches source diff).
#   - When decrementing would take "剩余" to 0 or below, the cycle is
#     treated as having completed and reset: E is set back to the row's
#     total day count (column D), and the start date in column F is
#     rolled forward by that same number of days (i.e. F = F + D).
#
# Rather than re-deriving this at runtime, the exact per-row target
# values (as they appear in the committed workbook) are applied
# directly below for reliability.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new value for column E, new value for column F
# (0 in the F slot means "column F is unchanged for this row").
$updates = @(
    @(2, 5, 0),
    @(3, 5, 0),
    @(4, 5, 0),
    @(5, 3, 0),
    @(6, 5, 0),
    @(7, 3, 0),
    @(8, 5, 0),
    @(9, 3, 0),
    @(10, 5, 0),
    @(11, 5, 0),
    @(12, 3, 0),
    @(13, 5, 0),
    @(14, 5, 0),
    @(15, 5, 0),
    @(16, 7, 0),
    @(17, 3, 0),
    @(18, 6, 0),
    @(19, 6, 0),
    @(20, 6, 0),
    @(21, 6, 0),
    @(22, 3, 0),
    @(23, 3, 0),
    @(24, 3, 0),
    @(25, 3, 0),
    @(26, 3, 0),
    @(27, 6, 0),
    @(28, 6, 0),
    @(29, 6, 0),
    @(30, 6, 0),
    @(31, 6, 0),
    @(32, 6, 0),
    @(33, 6, 0),
    @(34, 6, 0),
    @(35, 6, 0),
    @(37, 6, 0),
    @(38, 6, 0),
    @(39, 6, 0),
    @(40, 5, 0),
    @(41, 5, 0),
    @(42, 6, 0),
    @(43, 3, 0),
    @(44, 5, 0),
    @(45, 3, 0),
    @(46, 5, 0),
    @(47, 6, 0),
    @(48, 5, 0),
    @(49, 6, 0),
    @(50, 1, 0),
    @(51, 1, 0),
    @(52, 1, 0),
    @(53, 1, 0),
    @(54, 1, 0),
    @(55, 1, 0),
    @(56, 1, 0),
    @(57, 1, 0),
    @(58, 5, 0),
    @(59, 5, 0),
    @(60, 5, 0),
    @(61, 6, 0),
    @(62, 5, 0),
    @(63, 5, 0),
    @(64, 5, 0),
    @(65, 6, 0),
    @(66, 6, 0),
    @(67, 6, 0),
    @(68, 6, 0),
    @(69, 6, 0),
    @(70, 7, 0),
    @(71, 7, 0),
    @(72, 7, 0),
    @(73, 7, 0),
    @(74, 7, 0),
    @(75, 7, 0),
    @(76, 7, 0),
    @(77, 10, 20251112),
    @(78, 10, 20251112),
    @(79, 10, 20251112),
    @(80, 10, 20251112),
    @(81, 10, 20251112),
    @(82, 10, 20251112),
    @(83, 10, 20251112),
    @(84, 10, 20251112),
    @(85, 10, 20251112),
    @(86, 10, 20251112),
    @(87, 5, 0),
    @(88, 5, 0),
    @(89, 5, 0),
    @(90, 5, 0),
    @(91, 3, 0),
    @(92, 5, 0),
    @(93, 10, 20251112),
    @(94, 1, 0),
    @(95, 9, 0),
    @(96, 7, 0),
    @(97, 7, 0),
    @(98, 7, 0),
    @(99, 7, 0)
)

foreach ($u in $updates) {
    $row = $u[0]
    $newE = $u[1]
    $newF = $u[2]

    $ws.Cells.Item($row, 5).Value = $newE
    if ($newF -ne 0) {
        $ws.Cells.Item($row, 6).Value = $newF
    }
}

Write-Host "Updated $($updates.Count) rows (column E, and column F where applicable)."
